# Sample Project / Main.xlsx - "Rules" sheet
# The rule-name cell B11 changes from the text "R40" to the text "1".
#
# A plain  Range.Value = "1"  would be auto-coerced to the *number* 1 by
# Excel's type inference (since "1" parses as a number), which both changes
# the cell's stored type and forces a new cell style (quote-prefix / text
# number format) to be created. The original edit kept the existing cell
# style untouched and stored the new content as literal text, so we build
# the text value through a formula (which always yields a text result for
# a quoted string) and then paste it back as a value - this keeps the
# original number format/style on B11 and leaves the cell holding plain
# text "1", matching the recorded change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B11")

$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues - keep destination formatting/style
$excel.CutCopyMode = $false
